$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.574.48"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.645.53"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'213.06"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'0.533"
$ws.Range("E6").Value = "  +4.50%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'23.32"
$ws.Range("E8").Value = "  -4.08%  "
$ws.Range("D9").Value = "'0.256"
$ws.Range("E9").Value = "  -2.58%  "
$ws.Range("D10").Value = "'0.0611"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").Value = "'0.0889"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "1.879.66"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "1.642.66"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'0.585"
$ws.Range("E14").Value = "  +3.85%  "
$ws.Range("D15").Value = "'4.03"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "'64.41"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "27.552.26"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "'229.88"
$ws.Range("E18").Value = "  -4.38%  "
$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "'7.60"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'4.33"
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("D23").Value = "'9.72"
$ws.Range("E23").Value = "  +3.62%  "
$ws.Range("D24").Value = "'1.98"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").Value = "'148.86"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").Value = "'7.01"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'15.58"
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "'3.18"
$ws.Range("E33").Value = "  +3.14%  "
$ws.Range("D34").Value = "1.428.66"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").Value = "'0.571"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'0.884"
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("D39").Value = "'0.0167"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").Value = "'0.815"
$ws.Range("E42").Value = "  +3.38%  "
$ws.Range("D43").Value = "'5.47"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").Value = "'65.17"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("D46").Value = "1.788.81"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("D47").Value = "'1.68"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").Value = "'87.59"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").Value = "'0.0999"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").Value = "'7.75"
$ws.Range("E51").Value = "  -0.69%  "
